# Insert two new rows into the "Ají" price list right before the
# existing row 188, shifting all subsequent rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("188:189").Insert()

# Populate the newly inserted row 188 (a new "Primera" quality record
# for 2022-09-11 / serial 44806).
$ws.Cells.Item(188, 1).Value  = 8
$ws.Cells.Item(188, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(188, 3).Value  = "Coquimbo"
$ws.Cells.Item(188, 4).Value  = 44806
$ws.Cells.Item(188, 5).Value  = 4
$ws.Cells.Item(188, 6).Value  = 100112021
$ws.Cells.Item(188, 7).Value  = "Ají"
$ws.Cells.Item(188, 8).Value  = "Inferno"
$ws.Cells.Item(188, 9).Value  = "Primera"
$ws.Cells.Item(188, 10).Value = 500
$ws.Cells.Item(188, 11).Value = 16500
$ws.Cells.Item(188, 12).Value = 17000
$ws.Cells.Item(188, 13).Value = 16750
$ws.Cells.Item(188, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(188, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(188, 16).Value = 1675
$ws.Cells.Item(188, 17).Value = 10
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# Populate the newly inserted row 189 (a new "Segunda" quality record
# for the same date).
$ws.Cells.Item(189, 1).Value  = 8
$ws.Cells.Item(189, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(189, 3).Value  = "Coquimbo"
$ws.Cells.Item(189, 4).Value  = 44806
$ws.Cells.Item(189, 5).Value  = 4
$ws.Cells.Item(189, 6).Value  = 100112021
$ws.Cells.Item(189, 7).Value  = "Ají"
$ws.Cells.Item(189, 8).Value  = "Inferno"
$ws.Cells.Item(189, 9).Value  = "Segunda"
$ws.Cells.Item(189, 10).Value = 400
$ws.Cells.Item(189, 11).Value = 10500
$ws.Cells.Item(189, 12).Value = 11000
$ws.Cells.Item(189, 13).Value = 10750
$ws.Cells.Item(189, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(189, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(189, 16).Value = 1075
$ws.Cells.Item(189, 17).Value = 10
$ws.Cells.Item(189, 18).Value = "Hortaliza"

# Make sure the date cells keep the date/time number format used by
# the rest of column D.
$ws.Range("D188:D189").NumberFormat = $ws.Range("D187").NumberFormat
